$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 50,4
$data[0,0] = 0.06904570013284683
$data[0,1] = 0.9783502817153931
$data[0,2] = 0.02802884019911289
$data[0,3] = 0.9923202395439148
$data[1,0] = 0.00926543865352869
$data[1,1] = 0.9983460307121277
$data[1,2] = 0.01862300373613834
$data[1,3] = 0.9937091469764709
$data[2,0] = 0.005135770421475172
$data[2,1] = 0.9987913370132446
$data[2,2] = 0.007421276532113552
$data[2,3] = 0.9966503381729126
$data[3,0] = 0.002369061810895801
$data[3,1] = 0.9993850588798523
$data[3,2] = 0.003677255939692259
$data[3,3] = 0.998774528503418
$data[4,0] = 0.001487272442318499
$data[4,1] = 0.9997031092643738
$data[4,2] = 0.003892001695930958
$data[4,3] = 0.9988561868667603
$data[5,0] = 0.002277396852150559
$data[5,1] = 0.9995971322059631
$data[5,2] = 0.001003929879516363
$data[5,3] = 0.9996731877326965
$data[6,0] = 0.001231810194440186
$data[6,1] = 0.9996607303619385
$data[6,2] = 0.001089787343516946
$data[6,3] = 0.9997549057006836
$data[7,0] = 0.0009852364892140031
$data[7,1] = 0.9998091459274292
$data[7,2] = 0.006672353018075228
$data[7,3] = 0.9977124333381653
$data[8,0] = 0.001200216240249574
$data[8,1] = 0.9996607303619385
$data[8,2] = 0.003408713499084115
$data[8,3] = 0.998120903968811
$data[9,0] = 0.0006448489730246365
$data[9,1] = 0.9998303651809692
$data[9,2] = 0.0002576952974777669
$data[9,3] = 1
$data[10,0] = 0.0004412027483340353
$data[10,1] = 0.9999151825904846
$data[10,2] = 0.000298517057672143
$data[10,3] = 0.9998366236686707
$data[11,0] = 0.00103677180595696
$data[11,1] = 0.9997455477714539
$data[11,2] = 0.0007634095381945372
$data[11,3] = 0.9997549057006836
$data[12,0] = 0.000363870058208704
$data[12,1] = 0.9999151825904846
$data[12,2] = 0.0004264598537702113
$data[12,3] = 0.9999182820320129
$data[13,0] = 0.0006853825761936605
$data[13,1] = 0.9998727440834045
$data[13,2] = 0.0002707781677599996
$data[13,3] = 0.9999182820320129
$data[14,0] = 0.00006227634730748832
$data[14,1] = 1
$data[14,2] = 0.00141830206848681
$data[14,3] = 0.9993463754653931
$data[15,0] = 0.0006059093866497278
$data[15,1] = 0.9997879266738892
$data[15,2] = 0.001903250580653548
$data[15,3] = 0.9997549057006836
$data[16,0] = 0.0003233972238376737
$data[16,1] = 0.9998727440834045
$data[16,2] = 0.002583579858765006
$data[16,3] = 0.9993463754653931
$data[17,0] = 0.0004546408890746534
$data[17,1] = 0.9999151825904846
$data[17,2] = 0.003665063995867968
$data[17,3] = 0.9995915293693542
$data[18,0] = 0.0001014155786833726
$data[18,1] = 0.99997878074646
$data[18,2] = 0.0009293059119954705
$data[18,3] = 0.9999182820320129
$data[19,0] = 0.0003282561083324254
$data[19,1] = 0.9998939633369446
$data[19,2] = 0.0008562011062167585
$data[19,3] = 0.9997549057006836
$data[20,0] = 0.0003270724264439195
$data[20,1] = 0.9999364018440247
$data[20,2] = 0.002143299905583262
$data[20,3] = 0.9996731877326965
$data[21,0] = 0.00003890689549734816
$data[21,1] = 1
$data[21,2] = 0.001163584645837545
$data[21,3] = 0.9997549057006836
$data[22,0] = 0.000005231006980466191
$data[22,1] = 1
$data[22,2] = 0.001324944547377527
$data[22,3] = 0.9997549057006836
$data[23,0] = 0.0005088758189231157
$data[23,1] = 0.9998939633369446
$data[23,2] = 0.003270711982622743
$data[23,3] = 0.9997549057006836
$data[24,0] = 0.0001324149052379653
$data[24,1] = 0.9999151825904846
$data[24,2] = 0.001382868154905736
$data[24,3] = 0.9997549057006836
$data[25,0] = 0.0006425399915315211
$data[25,1] = 0.9998727440834045
$data[25,2] = 0.001903191790916026
$data[25,3] = 0.9997549057006836
$data[26,0] = 0.0001157905790023506
$data[26,1] = 0.9999364018440247
$data[26,2] = 0.002834173617884517
$data[26,3] = 0.9996731877326965
$data[27,0] = 0.0007897136965766549
$data[27,1] = 0.9998515844345093
$data[27,2] = 0.00139052071608603
$data[27,3] = 0.9997549057006836
$data[28,0] = 0.0002155366673832759
$data[28,1] = 0.9999364018440247
$data[28,2] = 0.0003383481525816023
$data[28,3] = 0.9998366236686707
$data[29,0] = 0.0002299517509527504
$data[29,1] = 0.9999364018440247
$data[29,2] = 0.0005281400517560542
$data[29,3] = 0.9997549057006836
$data[30,0] = 0.0001492913725087419
$data[30,1] = 0.9999364018440247
$data[30,2] = 0.001027223537676036
$data[30,3] = 0.9995915293693542
$data[31,0] = 0.00001127250470744912
$data[31,1] = 1
$data[31,2] = 0.0005563569138757885
$data[31,3] = 0.9998366236686707
$data[32,0] = 0.0007178789237514138
$data[32,1] = 0.9998727440834045
$data[32,2] = 0.0005911352345719934
$data[32,3] = 0.9998366236686707
$data[33,0] = 0.0003449516952969134
$data[33,1] = 0.9999575614929199
$data[33,2] = 0.0002298962790518999
$data[33,3] = 0.9999182820320129
$data[34,0] = 0.0001533441536594182
$data[34,1] = 0.9999575614929199
$data[34,2] = 0.0002323735825484619
$data[34,3] = 0.9999182820320129
$data[35,0] = 0.0001315911504207179
$data[35,1] = 0.9999575614929199
$data[35,2] = 0.0009848427725955844
$data[35,3] = 0.9997549057006836
$data[36,0] = 0.0007106171105988324
$data[36,1] = 0.9998727440834045
$data[36,2] = 0.0002411903842585161
$data[36,3] = 0.9999182820320129
$data[37,0] = 0.00009992806008085608
$data[37,1] = 0.99997878074646
$data[37,2] = 0.0001941750379046425
$data[37,3] = 0.9999182820320129
$data[38,0] = 0.000008794540008238982
$data[38,1] = 1
$data[38,2] = 0.00008571356011088938
$data[38,3] = 0.9999182820320129
$data[39,0] = 0.0002131737710442394
$data[39,1] = 0.9999364018440247
$data[39,2] = 0.001179247512482107
$data[39,3] = 0.9998366236686707
$data[40,0] = 0.0003750875184778124
$data[40,1] = 0.9998515844345093
$data[40,2] = 0.0006659092614427209
$data[40,3] = 0.9999182820320129
$data[41,0] = 0.0003728620358742774
$data[41,1] = 0.9999151825904846
$data[41,2] = 0.00002209134800068568
$data[41,3] = 1
$data[42,0] = 0.000006784814559068764
$data[42,1] = 1
$data[42,2] = 0.00005038882227381691
$data[42,3] = 1
$data[43,0] = 0.0001118665022659115
$data[43,1] = 0.9999575614929199
$data[43,2] = 0.0003329692408442497
$data[43,3] = 0.9998366236686707
$data[44,0] = 0.000005846857220603852
$data[44,1] = 1
$data[44,2] = 0.0009299021912738681
$data[44,3] = 0.9998366236686707
$data[45,0] = 0.00001133637761085993
$data[45,1] = 1
$data[45,2] = 0.0009652819717302918
$data[45,3] = 0.9998366236686707
$data[46,0] = 0.0002560740394983441
$data[46,1] = 0.9999151825904846
$data[46,2] = 0.001655063475482166
$data[46,3] = 0.9998366236686707
$data[47,0] = 0.0003139660984743387
$data[47,1] = 0.9999575614929199
$data[47,2] = 0.0008184186881408095
$data[47,3] = 0.9998366236686707
$data[48,0] = 0.0005276098381727934
$data[48,1] = 0.9999364018440247
$data[48,2] = 0.0009334432543255389
$data[48,3] = 0.9999182820320129
$data[49,0] = 0.0001398469612468034
$data[49,1] = 0.99997878074646
$data[49,2] = 0.0004749663057737052
$data[49,3] = 0.9999182820320129

$range = $ws.Range("A2:D51")
$range.Value = $data
